$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.245289444923401
$ws.Range("B1").Value = 2.321049213409424
$ws.Range("C1").Value = 3.084941625595093
$ws.Range("D1").Value = 3.559198617935181
$ws.Range("E1").Value = 1.379878282546997
